$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label text fixes (row 2) ---
$ws.Range("D2").Value = "Калории (к)"
$ws.Range("E2").Value = "Протеин (г)"

# --- Fill in missing numeric values for the first data row (row 3) ---
$ws.Range("D3").Value = 497
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 56

# --- Make the category header cells (row 1) bold ---
$ws.Range("D1").Font.Bold = $true
$ws.Range("H1").Font.Bold = $true
$ws.Range("Q1").Font.Bold = $true
$ws.Range("AF1").Font.Bold = $true
$ws.Range("AY1").Font.Bold = $true
$ws.Range("BD1").Font.Bold = $true
$ws.Range("BO1").Font.Bold = $true
$ws.Range("BT1").Font.Bold = $true

# --- Restore column E's width (side effect of the grouping below in the source edit) ---
$ws.Columns.Item(5).ColumnWidth = 12.1666667

# --- Group (outline) the detail columns D:BX so they can be collapsed/expanded ---
$ws.Range("D1:BX1").EntireColumn.Group()
